# Question 3 fix: the FY17/FY18/FY19 rank formulas (columns F, K, P) used a
# range reference that slid down with each row (e.g. E3:E53 on row 3) instead
# of staying pinned to the full data block (E2:E52). Anchor the ranges with
# absolute references so every row ranks against the same fixed range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# Row 2 holds its own (non-shared) formula; rows 3:52 form a shared-formula
# block. Setting them separately preserves that same grouping after the fix.

# FY17 rank (column F), ranked against FY17_diff_pct (column E)
$ws.Range("F2").Formula = '=RANK(E2,$E$2:$E$52,1)'
$ws.Range("F3:F52").Formula = '=RANK(E3,$E$2:$E$52,1)'

# FY18 rank (column K), ranked against FY18_diff_pct (column J)
$ws.Range("K2").Formula = '=RANK(J2,$J$2:$J$52,1)'
$ws.Range("K3:K52").Formula = '=RANK(J3,$J$2:$J$52,1)'

# FY19 rank (column P), ranked against FY19_diff_pct (column O)
$ws.Range("P2").Formula = '=RANK(O2,$O$2:$O$52,1)'
$ws.Range("P3:P52").Formula = '=RANK(O3,$O$2:$O$52,1)'

# Scroll the view down a bit (author was now working further down the sheet,
# around question 8) while keeping the existing selection (B92) untouched.
$ws.Activate()
$ws.Range("B92").Select()
$excel.ActiveWindow.ScrollRow = 73
$excel.ActiveWindow.ScrollColumn = 1
